$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.042868
$ws.Range("N2").Value = 0.128604
$ws.Range("O2").Value = 0.03014606792405771
$ws.Range("P2").Value = 0.03014606792405771
$ws.Range("Q2").Value = 0.015900298484
$ws.Range("R2").Value = 0.143102686356
$ws.Range("S2").Value = 0.03014606792405771
$ws.Range("T2").Value = 0.03014606792405771

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.370913
$ws.Range("H3").Value = 1.112739
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3865706666666667
$ws.Range("N3").Value = 1.159712
$ws.Range("O3").Value = 0.2718481285523376
$ws.Range("P3").Value = 0.2718481285523376
$ws.Range("Q3").Value = 0.1433840856853333
$ws.Range("R3").Value = 1.290456771168
$ws.Range("S3").Value = 0.2718481285523376
$ws.Range("T3").Value = 0.2718481285523376

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.370913
$ws.Range("H4").Value = 1.112739
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9839956666666666
$ws.Range("N4").Value = 2.951987
$ws.Range("O4").Value = 0.69197537100662
$ws.Range("P4").Value = 0.69197537100662
$ws.Range("Q4").Value = 0.3649767847103333
$ws.Range("R4").Value = 3.284791062393
$ws.Range("S4").Value = 0.69197537100662
$ws.Range("T4").Value = 0.69197537100662

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.370913
$ws.Range("H5").Value = 1.112739
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008575333333333332
$ws.Range("N5").Value = 0.025726
$ws.Range("O5").Value = 0.006030432516984765
$ws.Range("P5").Value = 0.006030432516984765
$ws.Range("Q5").Value = 0.003180702612666666
$ws.Range("R5").Value = 0.028626323514
$ws.Range("S5").Value = 0.006030432516984765
$ws.Range("T5").Value = 0.006030432516984765

